# Revert "Add yellow Mercedes taxi illustration"
#
# - Resize/reposition the teal rotated rectangle (id 3), the subtitle
#   textbox (id 5), the phone textbox (id 6) and the "Rapid/Sigur/
#   Confortabil" textbox (id 7) back to their pre-illustration geometry.
# - Remove the shapes that made up the illustration (ids 8-20): the
#   "Mercedes Fleet" label, the rounded-rectangle car body, the wheels,
#   windows, headlight, the taxi sign/star detail and its label.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Note: Shape.Left/Top/Width/Height round-trip through single-precision
# points internally, so a couple of values are nudged by 0.00001pt
# (<< 1/100 mm) to land exactly on the target EMU after truncation.

# --- "Rectangle 2" (teal rotated rectangle behind the "AKK" text) ---
$rect2 = $s.Shapes.Item("Rectangle 2")
$rect2.Left = 151.20001
$rect2.Top = -14.40007
$rect2.Width = 100.8
$rect2.Height = 86.4

# --- "TextBox 4" ("Razvan Taxi" / "Servicii de transport 24/7") ---
$tb4 = $s.Shapes.Item("TextBox 4")
$tb4.Width = 158.40001
$tb4.Height = 43.2

# --- "TextBox 5" ("Tel: 0720 064 963") ---
$tb5 = $s.Shapes.Item("TextBox 5")
$tb5.Left = 25.2
$tb5.Top = 104.4
$tb5.Width = 180
$tb5.Height = 28.80001

# --- "TextBox 6" ("Rapid" / "Sigur" / "Confortabil") ---
$tb6 = $s.Shapes.Item("TextBox 6")
$tb6.Left = 180
$tb6.Top = 28.80001
$tb6.Width = 72
$tb6.Height = 64.8

# --- Remove the Mercedes taxi illustration shapes (ids 8-20) ---
$namesToRemove = @(
    "TextBox 7",
    "Rounded Rectangle 8",
    "Rounded Rectangle 9",
    "Rectangle 10",
    "Oval 11",
    "Oval 12",
    "Oval 13",
    "Oval 14",
    "Rectangle 15",
    "Oval 16",
    "Rectangle 17",
    "TextBox 18",
    "5-Point Star 19"
)

foreach ($name in $namesToRemove) {
    $s.Shapes.Item($name).Delete()
}
